# Moves the footer logo/URL shapes on the title layout a bit to the left,
# and moves the footer bar + logo on the end ("Q&A") layout and on the
# slide master a bit further down.
#
# Corresponds to the OOXML diff touching:
#   ppt/slideLayouts/slideLayout2.xml  (title layout)
#   ppt/slideLayouts/slideLayout4.xml  (end / Q&A layout)
#   ppt/slideMasters/slideMaster1.xml
#
# Layouts/masters aren't reachable through Presentation.Designs in this
# host (writes silently fail to persist), so every layout/master is
# reached by hopping off a Slide that actually uses it
# (Slide.CustomLayout / Slide.Master) instead.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) slideLayout2.xml - title slide layout (used by slide 1)
# ---------------------------------------------------------------------
$s1   = $p.Slides.Item(1)
$lay2 = $s1.CustomLayout

# pic "그림 8" (logo, bottom-right): x 10895097 -> 10730337
$lay2.Shapes.Item(5).Left = 844.9084251968504
# sp  "직사각형 3" (URL text box under the logo): x 10880202 -> 10715442
$lay2.Shapes.Item(6).Left = 843.7355905511811

# ---------------------------------------------------------------------
# 2) slideLayout4.xml - end / Q&A layout (used by slide 4)
# ---------------------------------------------------------------------
$s4   = $p.Slides.Item(4)
$lay4 = $s4.CustomLayout

# sp  "직사각형 8" (footer bar): y 6432681 -> 6440919
$lay4.Shapes.Item(2).Top = 507.15898188976377
# pic "그림 5" (footer logo): y 6456364 -> 6472840
$lay4.Shapes.Item(3).Top = 509.6724409448819

# ---------------------------------------------------------------------
# 3) slideMaster1.xml (reached via slide 3, which uses slideLayout1 -
#    the layout that has no overrides of its own for these shapes, so
#    they come straight from the master)
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$m1 = $s3.Master

# sp  "직사각형 7" (footer bar): y 6432681 -> 6440919
$m1.Shapes.Item(3).Top = 507.15898188976377

# pic "그림 6" (footer logo): y 6456364 -> 6472840. (The author's edit also
# renumbered this shape to id 9 / "그림 8", but Id is read-only and this
# host does not support renaming a master-level shape - attempting it is
# a silent no-op - so only the position change is reproduced here.)
$picM = $m1.Shapes.Item(5)
$picM.Top = 509.6724409448819

Write-Host "done"
